$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("K2").Value = 1.8
$ws.Range("R2").Value = 1.3
$ws.Range("K3").Value = 1.73
$ws.Range("R3").Value = 1.25
$ws.Range("S3").Value = 1.85
$ws.Range("T3").Value = 1.95
$ws.Range("Q4").Value = 2.35
$ws.Range("R4").Value = 1.57
$ws.Range("G7").Value = 2
$ws.Range("I7").Value = 4.1
$ws.Range("S7").Value = 1.62
$ws.Range("T7").Value = 2.2
$ws.Range("Z7").Value = 17
$ws.Range("AA7").Value = 21
$ws.Range("AG7").Value = 8
$ws.Range("AH7").Value = 19
$ws.Range("AJ7").Value = 51
$ws.Range("AO7").Value = 12
$ws.Range("AT7").Value = 2.2
$ws.Range("AW7").Value = 6
$ws.Range("AX7").Value = 29
$ws.Range("G10").Value = 1.67
$ws.Range("H10").Value = 4.1
$ws.Range("I10").Value = 4.5
$ws.Range("J10").Value = 2.25
$ws.Range("K10").Value = 2.4
$ws.Range("L10").Value = 4.75
$ws.Range("O10").Value = 1.18
$ws.Range("P10").Value = 4.5
$ws.Range("Q10").Value = 1.62
$ws.Range("R10").Value = 2.25
$ws.Range("S10").Value = 1.27
$ws.Range("T10").Value = 3.4
$ws.Range("W10").Value = 9
$ws.Range("Z10").Value = 13
$ws.Range("AA10").Value = 12
$ws.Range("AD10").Value = 8
$ws.Range("AE10").Value = 13
$ws.Range("AH10").Value = 26
$ws.Range("AI10").Value = 15
$ws.Range("AJ10").Value = 51
$ws.Range("AO10").Value = 8.5
$ws.Range("AQ10").Value = 26
$ws.Range("AT10").Value = 3.4
$ws.Range("AV10").Value = 41
$ws.Range("BC10").Value = 451
$ws.Range("S11").Value = 1.3
